# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 8
    8  = 4
    9  = 6
    10 = 2
    11 = 11
    12 = 4
    13 = 2
    14 = 8
    15 = 4
    16 = 5
    17 = 4
    18 = 7
    19 = 4
    20 = 4
    21 = 5
    22 = 2
    23 = 8
    24 = 6
    25 = 1
    26 = 2
    27 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
